$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold plain text that looks numeric
# (e.g. "566.66", "1.955.59", "  +2.79%  "). Flip the range to Text
# format before writing so Excel stores the literal string instead of
# silently parsing it into a Double (which would corrupt values like
# "1.955.59" and introduce float rounding noise). Restore the style
# back to Normal afterwards so no visible formatting change remains.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.132.64"
$ws.Range("E2").Value = "  +2.79%  "

$ws.Range("D3").Value = "2.629.05"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "566.66"
$ws.Range("E5").Value = "  +6.01%  "

$ws.Range("D6").Value = "145.73"
$ws.Range("E6").Value = "  +2.59%  "

$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +6.84%  "

$ws.Range("D9").Value = "6.82"
$ws.Range("E9").Value = "  -2.05%  "

$ws.Range("D10").Value = "0.104"
$ws.Range("E10").Value = "  +3.32%  "

$ws.Range("E11").Value = "  +6.10%  "

$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").Value = "3.097.46"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").Value = "60.193.28"
$ws.Range("E14").Value = "  +3.01%  "

$ws.Range("D15").Value = "21.69"
$ws.Range("E15").Value = "  +4.07%  "

$ws.Range("D16").Value = "2.652.93"
$ws.Range("E16").Value = "  +1.56%  "

$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "4.56"
$ws.Range("E18").Value = "  +4.24%  "

$ws.Range("D19").Value = "343.71"
$ws.Range("E19").Value = "  +2.87%  "

$ws.Range("D20").Value = "10.37"
$ws.Range("E20").Value = "  +2.26%  "

$ws.Range("D21").Value = "6.27"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "66.79"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").Value = "0.435"
$ws.Range("E24").Value = "  +4.78%  "

$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").Value = "7.31"
$ws.Range("E27").Value = "  +2.90%  "

$ws.Range("D28").Value = "0.0₃0771"
$ws.Range("E28").Value = "  +5.43%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E30").Value = "  +4.23%  "

$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  +5.18%  "

$ws.Range("D32").Value = "157.89"
$ws.Range("E32").Value = "  +4.76%  "

$ws.Range("D33").Value = "19.11"
$ws.Range("E33").Value = "  +2.03%  "

$ws.Range("E34").Value = "  +4.76%  "

$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Value = "0.916"
$ws.Range("E35").Value = "  +11.74%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "0.912"
$ws.Range("E36").Value = "  +12.55%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.17"
$ws.Range("E37").Value = "  +6.65%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "1.52"
$ws.Range("E38").Value = "  +7.26%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "37.48"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").Value = "300.76"
$ws.Range("E40").Value = "  +6.64%  "

$ws.Range("D41").Value = "3.66"
$ws.Range("E41").Value = "  +2.32%  "

$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  -0.38%  "

$ws.Range("D43").Value = "0.604"
$ws.Range("E43").Value = "  +1.76%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.0974"
$ws.Range("E44").Value = "  +4.07%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0545"
$ws.Range("E45").Value = "  +2.84%  "

$ws.Range("D46").Value = "19.39"
$ws.Range("E46").Value = "  +2.24%  "

$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("D48").Value = "0.0236"
$ws.Range("E48").Value = "  +5.05%  "

$ws.Range("D49").Value = "4.71"
$ws.Range("E49").Value = "  +6.16%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "120.57"
$ws.Range("E50").Value = "  +7.82%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.955.59"
$ws.Range("E51").Value = "  +0.96%  "

$priceRange.Style = "Normal"
